$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> B
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9452054794520548
$ws.Range("C2").Value = 0.9583333333333334
$ws.Range("D2").Value = 0.9517241379310345
$ws.Range("E2").Value = 72

# Row 3 -> M
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.926829268292683
$ws.Range("C3").Value = 0.9047619047619048
$ws.Range("D3").Value = 0.9156626506024096
$ws.Range("E3").Value = 42

# Row 4 -> accuracy
$ws.Range("B4").Value = 0.9385964912280702
$ws.Range("C4").Value = 0.9385964912280702
$ws.Range("D4").Value = 0.9385964912280702
$ws.Range("E4").Value = 0.9385964912280702

# Row 5 -> macro avg
$ws.Range("B5").Value = 0.9360173738723689
$ws.Range("C5").Value = 0.9315476190476191
$ws.Range("D5").Value = 0.933693394266722
$ws.Range("E5").Value = 114

# Row 6 -> weighted avg
$ws.Range("B6").Value = 0.9384352963933389
$ws.Range("C6").Value = 0.9385964912280702
$ws.Range("D6").Value = 0.9384383268099622
$ws.Range("E6").Value = 114
